$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correcting data analysis: use the fixed initial %-control (0.3522) instead of
# chaining off the previous row's B-value for the normalisation formulas in
# column C (rows 19-27).
$ws.Range("C19").Formula = "=B19/0.3522*C18"
$ws.Range("C20:C27").Formula = "=B20/0.3522*C19"

# Daily entry: move the active selection/cursor to E21 and drop the stale
# scrolled-to-row16 viewport so the sheet opens at its natural top-left cell.
$ws.Range("E21").Select()
